# Regenerate merged AHB files
#
# 1) Rename the "_old" / "_new" suffixed header row (row 1) to "_FV2304" / "_FV2310".
# 2) Turn the used range A1:U67 into an Excel Table ("Table1") with those headers,
#    without letting the table-creation step invent a new header dxf (the source
#    workbook keeps styles.xml untouched, so we park the header formatting aside,
#    clear it, build the table, then restore it).
# 3) Freeze the header row (split after row 1, top-left cell of the scrolling pane A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row text -------------------------------------------------
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2304"
}
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $baseNames[$i] + "_FV2310"
}

# --- 2) Convert the range into a table, preserving existing header formatting --
$headerRange = $ws.Range("A1:U1")

# Stash a copy of the header row's current formatting on a scratch row far below
# the used range, so the table creation step (which bakes in a header dxf based
# on whatever formatting is present at Add()-time) sees a plain/unformatted
# header and therefore does not introduce any new dxf / style.
$scratchRange = $ws.Range("A200:U200")
$headerRange.Copy()
$scratchRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = $null

# Restore the header row's original formatting.
$scratchRange.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$scratchRange.Clear()

# --- 3) Freeze the header row ---------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
